$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 odds (columns G:BD)
$ws.Cells.Item(2, 7).Value = 1.65
$ws.Cells.Item(2, 8).Value = 3.9
$ws.Cells.Item(2, 9).Value = 4.7
$ws.Cells.Item(2, 10).Value = 2.18
$ws.Cells.Item(2, 11).Value = 2.25
$ws.Cells.Item(2, 12).Value = 4.85
$ws.Cells.Item(2, 13).Value = 1.05
$ws.Cells.Item(2, 14).Value = 8
$ws.Cells.Item(2, 15).Value = 1.25
$ws.Cells.Item(2, 16).Value = 3.55
$ws.Cells.Item(2, 17).Value = 1.78
$ws.Cells.Item(2, 18).Value = 1.98
$ws.Cells.Item(2, 19).Value = 1.37
$ws.Cells.Item(2, 20).Value = 2.87
$ws.Cells.Item(2, 21).Value = 1.78
$ws.Cells.Item(2, 22).Value = 1.93
$ws.Cells.Item(2, 23).Value = 7.4
$ws.Cells.Item(2, 24).Value = 7.9
$ws.Cells.Item(2, 25).Value = 8
$ws.Cells.Item(2, 26).Value = 12.5
$ws.Cells.Item(2, 27).Value = 12.5
$ws.Cells.Item(2, 28).Value = 25
$ws.Cells.Item(2, 29).Value = 8
$ws.Cells.Item(2, 30).Value = 7.6
$ws.Cells.Item(2, 31).Value = 16
$ws.Cells.Item(2, 32).Value = 70
$ws.Cells.Item(2, 33).Value = 500
$ws.Cells.Item(2, 34).Value = 13.5
$ws.Cells.Item(2, 35).Value = 27
$ws.Cells.Item(2, 36).Value = 15
$ws.Cells.Item(2, 37).Value = 80
$ws.Cells.Item(2, 38).Value = 45
$ws.Cells.Item(2, 39).Value = 45
$ws.Cells.Item(2, 40).Value = 3.5
$ws.Cells.Item(2, 41).Value = 7.9
$ws.Cells.Item(2, 42).Value = 17
$ws.Cells.Item(2, 43).Value = 25
$ws.Cells.Item(2, 44).Value = 55
$ws.Cells.Item(2, 45).Value = 250
$ws.Cells.Item(2, 46).Value = 2.87
$ws.Cells.Item(2, 47).Value = 7.6
$ws.Cells.Item(2, 48).Value = 70
$ws.Cells.Item(2, 49).Value = 6.4
$ws.Cells.Item(2, 50).Value = 26
$ws.Cells.Item(2, 51).Value = 32
$ws.Cells.Item(2, 52).Value = 150
$ws.Cells.Item(2, 53).Value = 175
$ws.Cells.Item(2, 54).Value = 400
$ws.Cells.Item(2, 55).Value = 81
$ws.Cells.Item(2, 56).Value = 81

# Update existing row 3 odds (columns G:BD)
$ws.Cells.Item(3, 7).Value = 1.65
$ws.Cells.Item(3, 8).Value = 3.4
$ws.Cells.Item(3, 9).Value = 6.25
$ws.Cells.Item(3, 10).Value = 2.3
$ws.Cells.Item(3, 11).Value = 2.05
$ws.Cells.Item(3, 12).Value = 6.5
$ws.Cells.Item(3, 13).Value = 1.1
$ws.Cells.Item(3, 14).Value = 7
$ws.Cells.Item(3, 15).Value = 1.44
$ws.Cells.Item(3, 16).Value = 2.63
$ws.Cells.Item(3, 17).Value = 2.35
$ws.Cells.Item(3, 18).Value = 1.57
$ws.Cells.Item(3, 19).Value = 1.53
$ws.Cells.Item(3, 20).Value = 2.38
$ws.Cells.Item(3, 21).Value = 2.25
$ws.Cells.Item(3, 22).Value = 1.57
$ws.Cells.Item(3, 23).Value = 5
$ws.Cells.Item(3, 24).Value = 6.5
$ws.Cells.Item(3, 25).Value = 9
$ws.Cells.Item(3, 26).Value = 12
$ws.Cells.Item(3, 27).Value = 17
$ws.Cells.Item(3, 28).Value = 41
$ws.Cells.Item(3, 29).Value = 7
$ws.Cells.Item(3, 30).Value = 7
$ws.Cells.Item(3, 31).Value = 21
$ws.Cells.Item(3, 32).Value = 81
$ws.Cells.Item(3, 33).Value = 101
$ws.Cells.Item(3, 34).Value = 12
$ws.Cells.Item(3, 35).Value = 29
$ws.Cells.Item(3, 36).Value = 21
$ws.Cells.Item(3, 37).Value = 67
$ws.Cells.Item(3, 38).Value = 51
$ws.Cells.Item(3, 39).Value = 67
$ws.Cells.Item(3, 40).Value = 3.4
$ws.Cells.Item(3, 41).Value = 9
$ws.Cells.Item(3, 42).Value = 26
$ws.Cells.Item(3, 43).Value = 29
$ws.Cells.Item(3, 44).Value = 67
$ws.Cells.Item(3, 45).Value = 251
$ws.Cells.Item(3, 46).Value = 2.38
$ws.Cells.Item(3, 47).Value = 10
$ws.Cells.Item(3, 48).Value = 81
$ws.Cells.Item(3, 49).Value = 7
$ws.Cells.Item(3, 50).Value = 34
$ws.Cells.Item(3, 51).Value = 41
$ws.Cells.Item(3, 52).Value = 151
$ws.Cells.Item(3, 53).Value = 201
$ws.Cells.Item(3, 54).Value = 51
$ws.Cells.Item(3, 55).Value = 51
$ws.Cells.Item(3, 56).Value = 51

# Insert two new rows before row 6 (shifts old rows 6,7 to 8,9)
$ws.Rows("6:7").Insert()

# Insert one new row after last data row (old row7 now at row9) to make room for row10
$ws.Rows("10:10").Insert()

# Fill new row 6 (lOwuTvTh - Thailand)
$ws.Cells.Item(6, 1).Value = "lOwuTvTh"
$ws.Cells.Item(6, 2).Value = "24/11/2024"
$ws.Cells.Item(6, 3).Value = "08:00"
$ws.Cells.Item(6, 4).Value = "THAILAND - THAI LEAGUE 1"
$ws.Cells.Item(6, 5).Value = "Nong Bua Pitchaya"
$ws.Cells.Item(6, 6).Value = "Chiangrai Utd"
$ws.Cells.Item(6, 7).Value = 1.8
$ws.Cells.Item(6, 8).Value = 3.75
$ws.Cells.Item(6, 9).Value = 3.85
$ws.Cells.Item(6, 10).Value = 2.32
$ws.Cells.Item(6, 11).Value = 2.3
$ws.Cells.Item(6, 12).Value = 4.1
$ws.Cells.Item(6, 13).Value = 1.04
$ws.Cells.Item(6, 14).Value = 8.5
$ws.Cells.Item(6, 15).Value = 1.21
$ws.Cells.Item(6, 16).Value = 4
$ws.Cells.Item(6, 17).Value = 1.62
$ws.Cells.Item(6, 18).Value = 2.2
$ws.Cells.Item(6, 19).Value = 1.32
$ws.Cells.Item(6, 20).Value = 3.1
$ws.Cells.Item(6, 21).Value = 1.6
$ws.Cells.Item(6, 22).Value = 2.22
$ws.Cells.Item(6, 23).Value = 9
$ws.Cells.Item(6, 24).Value = 10
$ws.Cells.Item(6, 25).Value = 8.25
$ws.Cells.Item(6, 26).Value = 15.5
$ws.Cells.Item(6, 27).Value = 13
$ws.Cells.Item(6, 28).Value = 20
$ws.Cells.Item(6, 29).Value = 8.5
$ws.Cells.Item(6, 30).Value = 7.4
$ws.Cells.Item(6, 31).Value = 13
$ws.Cells.Item(6, 32).Value = 45
$ws.Cells.Item(6, 33).Value = 300
$ws.Cells.Item(6, 34).Value = 14
$ws.Cells.Item(6, 35).Value = 23
$ws.Cells.Item(6, 36).Value = 13
$ws.Cells.Item(6, 37).Value = 55
$ws.Cells.Item(6, 38).Value = 32
$ws.Cells.Item(6, 39).Value = 32
$ws.Cells.Item(6, 40).Value = 3.9
$ws.Cells.Item(6, 41).Value = 8.75
$ws.Cells.Item(6, 42).Value = 15.5
$ws.Cells.Item(6, 43).Value = 29
$ws.Cells.Item(6, 44).Value = 50
$ws.Cells.Item(6, 45).Value = 175
$ws.Cells.Item(6, 46).Value = 3.1
$ws.Cells.Item(6, 47).Value = 6.8
$ws.Cells.Item(6, 48).Value = 50
$ws.Cells.Item(6, 49).Value = 5.9
$ws.Cells.Item(6, 50).Value = 21
$ws.Cells.Item(6, 51).Value = 24
$ws.Cells.Item(6, 52).Value = 100
$ws.Cells.Item(6, 53).Value = 120
$ws.Cells.Item(6, 54).Value = 250
$ws.Cells.Item(6, 55).Value = ""
$ws.Cells.Item(6, 56).Value = ""

# Fill new row 7 (I3VePdSH - Thailand)
$ws.Cells.Item(7, 1).Value = "I3VePdSH"
$ws.Cells.Item(7, 2).Value = "24/11/2024"
$ws.Cells.Item(7, 3).Value = "08:00"
$ws.Cells.Item(7, 4).Value = "THAILAND - THAI LEAGUE 1"
$ws.Cells.Item(7, 5).Value = "Ratchaburi"
$ws.Cells.Item(7, 6).Value = "Khonkaen Utd."
$ws.Cells.Item(7, 7).Value = 1.37
$ws.Cells.Item(7, 8).Value = 4.65
$ws.Cells.Item(7, 9).Value = 7.4
$ws.Cells.Item(7, 10).Value = 1.83
$ws.Cells.Item(7, 11).Value = 2.45
$ws.Cells.Item(7, 12).Value = 6.5
$ws.Cells.Item(7, 13).Value = 1.04
$ws.Cells.Item(7, 14).Value = 8.75
$ws.Cells.Item(7, 15).Value = 1.19
$ws.Cells.Item(7, 16).Value = 4.15
$ws.Cells.Item(7, 17).Value = 1.6
$ws.Cells.Item(7, 18).Value = 2.22
$ws.Cells.Item(7, 19).Value = 1.31
$ws.Cells.Item(7, 20).Value = 3.15
$ws.Cells.Item(7, 21).Value = 1.83
$ws.Cells.Item(7, 22).Value = 1.87
$ws.Cells.Item(7, 23).Value = 7.8
$ws.Cells.Item(7, 24).Value = 7
$ws.Cells.Item(7, 25).Value = 8.25
$ws.Cells.Item(7, 26).Value = 9
$ws.Cells.Item(7, 27).Value = 10.75
$ws.Cells.Item(7, 28).Value = 24
$ws.Cells.Item(7, 29).Value = 8.75
$ws.Cells.Item(7, 30).Value = 9.25
$ws.Cells.Item(7, 31).Value = 18.5
$ws.Cells.Item(7, 32).Value = 80
$ws.Cells.Item(7, 33).Value = 600
$ws.Cells.Item(7, 34).Value = 21
$ws.Cells.Item(7, 35).Value = 50
$ws.Cells.Item(7, 36).Value = 23
$ws.Cells.Item(7, 37).Value = 175
$ws.Cells.Item(7, 38).Value = 80
$ws.Cells.Item(7, 39).Value = 65
$ws.Cells.Item(7, 40).Value = 3.3
$ws.Cells.Item(7, 41).Value = 6.2
$ws.Cells.Item(7, 42).Value = 15.5
$ws.Cells.Item(7, 43).Value = 16.5
$ws.Cells.Item(7, 44).Value = 40
$ws.Cells.Item(7, 45).Value = 200
$ws.Cells.Item(7, 46).Value = 3.15
$ws.Cells.Item(7, 47).Value = 8
$ws.Cells.Item(7, 48).Value = 70
$ws.Cells.Item(7, 49).Value = 8.5
$ws.Cells.Item(7, 50).Value = 40
$ws.Cells.Item(7, 51).Value = 40
$ws.Cells.Item(7, 52).Value = 300
$ws.Cells.Item(7, 53).Value = 250
$ws.Cells.Item(7, 54).Value = 500
$ws.Cells.Item(7, 55).Value = ""
$ws.Cells.Item(7, 56).Value = ""

# Fill row 8 (moved/modified xppFPDhg - Turkey Super Lig)
$ws.Cells.Item(8, 1).Value = "xppFPDhg"
$ws.Cells.Item(8, 2).Value = "24/11/2024"
$ws.Cells.Item(8, 3).Value = "07:30"
$ws.Cells.Item(8, 4).Value = "TURKEY - SUPER LIG"
$ws.Cells.Item(8, 5).Value = "Sivasspor"
$ws.Cells.Item(8, 6).Value = "Kasimpasa"
$ws.Cells.Item(8, 7).Value = 3.1
$ws.Cells.Item(8, 8).Value = 3.5
$ws.Cells.Item(8, 9).Value = 2.2
$ws.Cells.Item(8, 10).Value = 3.6
$ws.Cells.Item(8, 11).Value = 2.25
$ws.Cells.Item(8, 12).Value = 2.88
$ws.Cells.Item(8, 13).Value = 1.04
$ws.Cells.Item(8, 14).Value = 13
$ws.Cells.Item(8, 15).Value = 1.22
$ws.Cells.Item(8, 16).Value = 4
$ws.Cells.Item(8, 17).Value = 1.75
$ws.Cells.Item(8, 18).Value = 2.05
$ws.Cells.Item(8, 19).Value = 1.36
$ws.Cells.Item(8, 20).Value = 3
$ws.Cells.Item(8, 21).Value = 1.62
$ws.Cells.Item(8, 22).Value = 2.2
$ws.Cells.Item(8, 23).Value = 11
$ws.Cells.Item(8, 24).Value = 17
$ws.Cells.Item(8, 25).Value = 11
$ws.Cells.Item(8, 26).Value = 34
$ws.Cells.Item(8, 27).Value = 23
$ws.Cells.Item(8, 28).Value = 29
$ws.Cells.Item(8, 29).Value = 12
$ws.Cells.Item(8, 30).Value = 6.5
$ws.Cells.Item(8, 31).Value = 13
$ws.Cells.Item(8, 32).Value = 41
$ws.Cells.Item(8, 33).Value = 151
$ws.Cells.Item(8, 34).Value = 9.5
$ws.Cells.Item(8, 35).Value = 12
$ws.Cells.Item(8, 36).Value = 9
$ws.Cells.Item(8, 37).Value = 21
$ws.Cells.Item(8, 38).Value = 17
$ws.Cells.Item(8, 39).Value = 23
$ws.Cells.Item(8, 40).Value = 5
$ws.Cells.Item(8, 41).Value = 17
$ws.Cells.Item(8, 42).Value = 23
$ws.Cells.Item(8, 43).Value = 51
$ws.Cells.Item(8, 44).Value = 67
$ws.Cells.Item(8, 45).Value = 151
$ws.Cells.Item(8, 46).Value = 3
$ws.Cells.Item(8, 47).Value = 7.5
$ws.Cells.Item(8, 48).Value = 51
$ws.Cells.Item(8, 49).Value = 4.33
$ws.Cells.Item(8, 50).Value = 12
$ws.Cells.Item(8, 51).Value = 21
$ws.Cells.Item(8, 52).Value = 41
$ws.Cells.Item(8, 53).Value = 51
$ws.Cells.Item(8, 54).Value = 126
$ws.Cells.Item(8, 55).Value = 251
$ws.Cells.Item(8, 56).Value = 301

# Fill row 9 (moved/modified d2jjMXa3 - Turkey 1. Lig)
$ws.Cells.Item(9, 1).Value = "d2jjMXa3"
$ws.Cells.Item(9, 2).Value = "24/11/2024"
$ws.Cells.Item(9, 3).Value = "07:30"
$ws.Cells.Item(9, 4).Value = "TURKEY - 1. LIG"
$ws.Cells.Item(9, 5).Value = "Erzurumspor"
$ws.Cells.Item(9, 6).Value = "Karagumruk"
$ws.Cells.Item(9, 7).Value = 2.2
$ws.Cells.Item(9, 8).Value = 3.3
$ws.Cells.Item(9, 9).Value = 3.1
$ws.Cells.Item(9, 10).Value = 3
$ws.Cells.Item(9, 11).Value = 2.05
$ws.Cells.Item(9, 12).Value = 3.75
$ws.Cells.Item(9, 13).Value = 1.07
$ws.Cells.Item(9, 14).Value = 9
$ws.Cells.Item(9, 15).Value = 1.33
$ws.Cells.Item(9, 16).Value = 3.25
$ws.Cells.Item(9, 17).Value = 2.1
$ws.Cells.Item(9, 18).Value = 1.7
$ws.Cells.Item(9, 19).Value = 1.44
$ws.Cells.Item(9, 20).Value = 2.63
$ws.Cells.Item(9, 21).Value = 1.83
$ws.Cells.Item(9, 22).Value = 1.83
$ws.Cells.Item(9, 23).Value = 7
$ws.Cells.Item(9, 24).Value = 10
$ws.Cells.Item(9, 25).Value = 9.5
$ws.Cells.Item(9, 26).Value = 21
$ws.Cells.Item(9, 27).Value = 19
$ws.Cells.Item(9, 28).Value = 29
$ws.Cells.Item(9, 29).Value = 9
$ws.Cells.Item(9, 30).Value = 6.5
$ws.Cells.Item(9, 31).Value = 15
$ws.Cells.Item(9, 32).Value = 51
$ws.Cells.Item(9, 33).Value = 301
$ws.Cells.Item(9, 34).Value = 9
$ws.Cells.Item(9, 35).Value = 15
$ws.Cells.Item(9, 36).Value = 12
$ws.Cells.Item(9, 37).Value = 34
$ws.Cells.Item(9, 38).Value = 26
$ws.Cells.Item(9, 39).Value = 34
$ws.Cells.Item(9, 40).Value = 4.33
$ws.Cells.Item(9, 41).Value = 13
$ws.Cells.Item(9, 42).Value = 23
$ws.Cells.Item(9, 43).Value = 41
$ws.Cells.Item(9, 44).Value = 67
$ws.Cells.Item(9, 45).Value = 201
$ws.Cells.Item(9, 46).Value = 2.63
$ws.Cells.Item(9, 47).Value = 8
$ws.Cells.Item(9, 48).Value = 51
$ws.Cells.Item(9, 49).Value = 5
$ws.Cells.Item(9, 50).Value = 19
$ws.Cells.Item(9, 51).Value = 29
$ws.Cells.Item(9, 52).Value = 51
$ws.Cells.Item(9, 53).Value = 81
$ws.Cells.Item(9, 54).Value = 201
$ws.Cells.Item(9, 55).Value = 126
$ws.Cells.Item(9, 56).Value = 126

# Fill new row 10 (2yBJ1dI0 - Ukraine Premier League)
$ws.Cells.Item(10, 1).Value = "2yBJ1dI0"
$ws.Cells.Item(10, 2).Value = "24/11/2024"
$ws.Cells.Item(10, 3).Value = "08:00"
$ws.Cells.Item(10, 4).Value = "UKRAINE - PREMIER LEAGUE"
$ws.Cells.Item(10, 5).Value = "Obolon"
$ws.Cells.Item(10, 6).Value = "Kryvbas"
$ws.Cells.Item(10, 7).Value = 4.5
$ws.Cells.Item(10, 8).Value = 3.4
$ws.Cells.Item(10, 9).Value = 1.75
$ws.Cells.Item(10, 10).Value = 5
$ws.Cells.Item(10, 11).Value = 2.02
$ws.Cells.Item(10, 12).Value = 2.37
$ws.Cells.Item(10, 13).Value = 1.08
$ws.Cells.Item(10, 14).Value = 7.85
$ws.Cells.Item(10, 15).Value = 1.42
$ws.Cells.Item(10, 16).Value = 2.47
$ws.Cells.Item(10, 17).Value = 2.22
$ws.Cells.Item(10, 18).Value = 1.52
$ws.Cells.Item(10, 19).Value = 1.47
$ws.Cells.Item(10, 20).Value = 2.32
$ws.Cells.Item(10, 21).Value = 2.12
$ws.Cells.Item(10, 22).Value = 1.57
$ws.Cells.Item(10, 23).Value = 9.75
$ws.Cells.Item(10, 24).Value = 23
$ws.Cells.Item(10, 25).Value = 16
$ws.Cells.Item(10, 26).Value = 80
$ws.Cells.Item(10, 27).Value = 55
$ws.Cells.Item(10, 28).Value = 75
$ws.Cells.Item(10, 29).Value = 7.4
$ws.Cells.Item(10, 30).Value = 6.8
$ws.Cells.Item(10, 31).Value = 21
$ws.Cells.Item(10, 32).Value = 150
$ws.Cells.Item(10, 33).Value = 900
$ws.Cells.Item(10, 34).Value = 5.3
$ws.Cells.Item(10, 35).Value = 6.9
$ws.Cells.Item(10, 36).Value = 8.75
$ws.Cells.Item(10, 37).Value = 13
$ws.Cells.Item(10, 38).Value = 16.5
$ws.Cells.Item(10, 39).Value = 40
$ws.Cells.Item(10, 40).Value = 6
$ws.Cells.Item(10, 41).Value = 27
$ws.Cells.Item(10, 42).Value = 40
$ws.Cells.Item(10, 43).Value = 175
$ws.Cells.Item(10, 44).Value = 250
$ws.Cells.Item(10, 45).Value = 500
$ws.Cells.Item(10, 46).Value = 2.3
$ws.Cells.Item(10, 47).Value = 8.5
$ws.Cells.Item(10, 48).Value = 100
$ws.Cells.Item(10, 49).Value = 3.35
$ws.Cells.Item(10, 50).Value = 8.75
$ws.Cells.Item(10, 51).Value = 22
$ws.Cells.Item(10, 52).Value = 32
$ws.Cells.Item(10, 53).Value = 80
$ws.Cells.Item(10, 54).Value = 400
$ws.Cells.Item(10, 55).Value = 81
$ws.Cells.Item(10, 56).Value = 81

